$d = $word.ActiveDocument

# Update the date heading in the first paragraph.
$dateRange = $d.Paragraphs.Item(1).Range
$dateRange.MoveEnd(1, -1) | Out-Null
$dateRange.Text = "2025-02-11 Tuesday"

# Update each math-problem cell in the table, addressed by (row, col)
# rather than text search, since several problems repeat their operands
# (e.g. two different cells both originally read "26+24=") and a plain
# Find/Replace would not be able to tell them apart.
$t = $d.Tables.Item(1)

$cellRange = $t.Cell(1,1).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "64-56="

$cellRange = $t.Cell(1,2).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "49+29="

$cellRange = $t.Cell(1,3).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "0+29="

$cellRange = $t.Cell(1,4).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "49+23="

$cellRange = $t.Cell(1,5).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "47-44="

$cellRange = $t.Cell(2,1).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "28-8="

$cellRange = $t.Cell(2,2).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "14+79="

$cellRange = $t.Cell(2,3).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "75-48="

$cellRange = $t.Cell(2,4).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "75-3="

$cellRange = $t.Cell(2,5).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "30-11="

$cellRange = $t.Cell(3,1).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "4+86="

$cellRange = $t.Cell(3,2).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "47+21="

$cellRange = $t.Cell(3,3).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "93-49="

$cellRange = $t.Cell(3,4).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "58-2="

$cellRange = $t.Cell(3,5).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "16+8="

$cellRange = $t.Cell(4,1).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "88-66="

$cellRange = $t.Cell(4,2).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "40+3="

$cellRange = $t.Cell(4,3).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "57+9="

$cellRange = $t.Cell(4,4).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "5+89="

$cellRange = $t.Cell(4,5).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "80-1="

$cellRange = $t.Cell(5,1).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "89+7="

$cellRange = $t.Cell(5,2).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "77+9="

$cellRange = $t.Cell(5,3).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "83-49="

$cellRange = $t.Cell(5,4).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "18-18="

$cellRange = $t.Cell(5,5).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "75-51="

$cellRange = $t.Cell(6,1).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "92-43="

$cellRange = $t.Cell(6,2).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "68-35="

$cellRange = $t.Cell(6,3).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "57-46="

$cellRange = $t.Cell(6,4).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "15+55="

$cellRange = $t.Cell(6,5).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "23-14="

$cellRange = $t.Cell(7,1).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "58-16="

$cellRange = $t.Cell(7,2).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "54-7="

$cellRange = $t.Cell(7,3).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "27+20="

$cellRange = $t.Cell(7,4).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "94-64="

$cellRange = $t.Cell(7,5).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "74-23="

$cellRange = $t.Cell(8,1).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "65-15="

$cellRange = $t.Cell(8,2).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "46-7="

$cellRange = $t.Cell(8,3).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "98-27="

$cellRange = $t.Cell(8,4).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "68+24="

$cellRange = $t.Cell(8,5).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "44+41="

$cellRange = $t.Cell(9,1).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "90-82="

$cellRange = $t.Cell(9,2).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "55-4="

$cellRange = $t.Cell(9,3).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "60+19="

$cellRange = $t.Cell(9,4).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "19+52="

$cellRange = $t.Cell(9,5).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "2+81="

$cellRange = $t.Cell(10,1).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "51-35="

$cellRange = $t.Cell(10,2).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "65+13="

$cellRange = $t.Cell(10,3).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "68+13="

$cellRange = $t.Cell(10,4).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "86-46="

$cellRange = $t.Cell(10,5).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "86-69="

$cellRange = $t.Cell(11,1).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "71-52="

$cellRange = $t.Cell(11,2).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "67+6="

$cellRange = $t.Cell(11,3).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "24-3="

$cellRange = $t.Cell(11,4).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "13+64="

$cellRange = $t.Cell(11,5).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "11+88="

$cellRange = $t.Cell(12,1).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "48-7="

$cellRange = $t.Cell(12,2).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "57-44="

$cellRange = $t.Cell(12,3).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "11+60="

$cellRange = $t.Cell(12,4).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "65+19="

$cellRange = $t.Cell(12,5).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "85-27="

$cellRange = $t.Cell(13,1).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "26-16="

$cellRange = $t.Cell(13,2).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "72-20="

$cellRange = $t.Cell(13,3).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "99-29="

$cellRange = $t.Cell(13,4).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "49-38="

$cellRange = $t.Cell(13,5).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "28+37="

$cellRange = $t.Cell(14,1).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "76-6="

$cellRange = $t.Cell(14,2).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "31+46="

$cellRange = $t.Cell(14,3).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "49-9="

$cellRange = $t.Cell(14,4).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "57+33="

$cellRange = $t.Cell(14,5).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "97-72="

$cellRange = $t.Cell(15,1).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "28+1="

$cellRange = $t.Cell(15,2).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "19+12="

$cellRange = $t.Cell(15,3).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "30-25="

$cellRange = $t.Cell(15,4).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "80-40="

$cellRange = $t.Cell(15,5).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "34-33="

$cellRange = $t.Cell(16,1).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "12+59="

$cellRange = $t.Cell(16,2).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "94-28="

$cellRange = $t.Cell(16,3).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "66-15="

$cellRange = $t.Cell(16,4).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "37-8="

$cellRange = $t.Cell(16,5).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "93-42="

$cellRange = $t.Cell(17,1).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "47+33="

$cellRange = $t.Cell(17,2).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "6+30="

$cellRange = $t.Cell(17,3).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "84+11="

$cellRange = $t.Cell(17,4).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "1+42="

$cellRange = $t.Cell(17,5).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "86-33="

$cellRange = $t.Cell(18,1).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "33-22="

$cellRange = $t.Cell(18,2).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "16+20="

$cellRange = $t.Cell(18,3).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "54-10="

$cellRange = $t.Cell(18,4).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "39-11="

$cellRange = $t.Cell(18,5).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "68+24="

$cellRange = $t.Cell(19,1).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "46+16="

$cellRange = $t.Cell(19,2).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "40-24="

$cellRange = $t.Cell(19,3).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "3+81="

$cellRange = $t.Cell(19,4).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "30+51="

$cellRange = $t.Cell(19,5).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "40-12="

$cellRange = $t.Cell(20,1).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "5+38="

$cellRange = $t.Cell(20,2).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "1+44="

$cellRange = $t.Cell(20,3).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "36+54="

$cellRange = $t.Cell(20,4).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "65-47="

$cellRange = $t.Cell(20,5).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "15+50="
